$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Fiber Type"): rows 4 and 5 change from "x" to "x+x"
$ws.Range("E4").Value = "x+x"
$ws.Range("E5").Value = "x+x"

# Column G ("Beta"): numeric values become text "<val>+<val>"
$ws.Range("G2").Value = "2+2"
$ws.Range("G3").Value = "3.3+3.3"
$ws.Range("G4").Value = "9.1+9.1"
$ws.Range("G5").Value = "2+2"

# Column H ("Gamma"): numeric values become text "<val>+<val>"
$ws.Range("H2").Value = "0.25+0.25"
$ws.Range("H3").Value = "0.63+0.63"
$ws.Range("H4").Value = "8.1+8.1"
$ws.Range("H5").Value = "0.14+0.14"

# Column I ("Dispersion"): numeric 1 becomes text "1+1"
$ws.Range("I2").Value = "1+1"
$ws.Range("I3").Value = "1+1"
$ws.Range("I4").Value = "1+1"
$ws.Range("I5").Value = "1+1"

# Update the selected cell shown in the sheet view
$ws.Range("G8").Select()
